$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.217.34"
$ws.Range("E2").Value = "  +0.93%  "
$ws.Range("D3").Value = "1.852.70"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.16"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4625"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3718"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8875"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "1.982.17"
$ws.Range("E11").Value = "  +8.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.01"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07801"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.371"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.510"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.17"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008902"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.70"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "27.242.93"
$ws.Range("E21").Value = "  +0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.058"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.50"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "2.065.98"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.948"
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.10"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.046"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.70"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.058"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08825"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.125"
$ws.Range("E32").Value = "  +5.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7646"
$ws.Range("E33").Value = "  +4.96%  "
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.498"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.751"
$ws.Range("E36").Value = "  +10.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05237"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.052"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1625"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.396"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.34"
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.14"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.637"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06201"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.76"
$ws.Range("E51").Value = "  +1.46%  "
